$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Work from the bottom of the sheet upward so that row numbers used below
# (which refer to the ORIGINAL layout) stay valid while we edit.

# --- Area around original row 22 (Noble Park) ---
# Remove the "Noble Park" row and replace it with two new "Pascoe Vale" rows.
$ws.Rows.Item(22).Delete()

$ws.Rows.Item(22).Insert()
$ws.Rows.Item(22).Insert()

$ws.Range("A22").Value = "Pascoe Vale"
$ws.Range("B22").Value = "Elite Swimming Pascoe Vale, 8 Attercliffe Avenue"
$ws.Range("C22").Value = "5pm - 6pm 8/2/2021"
$ws.Range("D22").Value = "Case attended venue"

$ws.Range("A23").Value = "Pascoe Vale"
$ws.Range("B23").Value = "Oak Park Sports and Aquatic Centre, 563a Pascoe Vale Road"
$ws.Range("C23").Value = "4pm - 7.30pm 10/2/2021"
$ws.Range("D23").Value = "Case attended venue"

# --- Original row 15 (Keysborough - Aces Sporting Club) ---
# This whole venue/time entry is removed from the table.
$ws.Rows.Item(15).Delete()

# --- Area around original row 7 (right after Broadmeadows train row) ---
# Insert two new "Broadmeadows" venue rows.
$ws.Rows.Item(7).Insert()
$ws.Rows.Item(7).Insert()

$ws.Range("A7").Value = "Broadmeadows"
$ws.Range("B7").Value = "Ferguson Plarre Bakehouses - Broadmeadows, 1099-1169 Pascoe Vale Road"
$ws.Range("C7").Value = "12:30pm - 12:45pm 9/2/2021"
$ws.Range("D7").Value = "Case attended venue"

$ws.Range("A8").Value = "Broadmeadows"
$ws.Range("B8").Value = "Woolworths Broadmeadows Central, Pascoe Vale Road"
$ws.Range("C8").Value = "12.15pm - 12:30 pm 9/2/2021"
$ws.Range("D8").Value = "Case attended venue"

# --- Fix typo in exposure period for row 6 (Broadmeadows / Craigieburn Line train) ---
$ws.Range("C6").Value = "1.25pm - 1.59pm  9/02/2021"
